$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka1")

# Update the TODO column (C) with the reshuffled task list for weeks 2-14 (rows 3-15)
$ws.Range("C3").Value  = "Termék filterezés, landing page, UI szépítése"
$ws.Range("C4").Value  = "Deployment, hosting, webes fizetés"
$ws.Range("C5").Value  = "Deployment, hosting, webes fizetés, tesztek készítése"
$ws.Range("C6").Value  = "Webes fizetés, tesztek készítése"
$ws.Range("C7").Value  = "Android kliens - architektura megtervezés, app skeleton létrehozása"
$ws.Range("C8").Value  = "Android kliens - kezdőképernyő, autentikáció megvalósítás, lokális adatbázis megvalósítás"
$ws.Range("C9").Value  = "Android kliens - termékek, kosár, checkout képernyők megvalósítása"
$ws.Range("C10").Value = "Android kliens - termék filterezés megvalósítása"
$ws.Range("C11").Value = "Diplomamunka írás"
$ws.Range("C12").Value = "Diplomamunka írás"
$ws.Range("C13").Value = "Diplomamunka írás"
$ws.Range("C14").Value = "Diplomamunka írás"
$ws.Range("C15").Value = "Diplomamunka írás"

# Update the active selection to C6, matching the saved view state
$ws.Range("C6").Select()
